# Apply the "Pong Game" test-case section to the GameTestCases workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Select Level section: row 11's Sprint # changes from 2 to 1
# ---------------------------------------------------------------------
$ws.Range("A11").Value = 1

# ---------------------------------------------------------------------
# 2. Snake MG section: fill in the (previously empty) Pass/Fail column
#    for several rows with "pass"
# ---------------------------------------------------------------------
$ws.Range("E22").Value = "pass"
$ws.Range("E23").Value = "pass"
$ws.Range("E27").Value = "pass"
$ws.Range("E28").Value = "pass"
$ws.Range("E29").Value = "pass"
$ws.Range("E30").Value = "pass"
$ws.Range("E31").Value = "pass"
$ws.Range("E32").Value = "pass"

# ---------------------------------------------------------------------
# 3. New "Pong Game" test-case block, rows 34-45.
#    Cell values are entered in the same order the original author used
#    so that newly-created shared strings line up the same way.
# ---------------------------------------------------------------------
$ws.Range("D34").Value = "ball starts moving"
$ws.Range("D35").Value = "right bar moves up"
$ws.Range("D36").Value = "right bar moves down"
$ws.Range("C37").Value = "w key"
$ws.Range("D37").Value = "left bar moves up"
$ws.Range("C38").Value = "s key"
$ws.Range("D38").Value = "left bar moves down"
$ws.Range("D39").Value = "ball moves in opposite direction"
$ws.Range("D40").Value = "displays game over with winner player "
$ws.Range("C44").Value = "Main Menu Button Hover"
$ws.Range("C45").Value = "Main Menu Button Press"
$ws.Range("D45").Value = "highlights grey on mouse hold and redirects player to main menu"
$ws.Range("C34").Value = "Player waits 2.5 seconds"
$ws.Range("C35").Value = "Up arrow key"
$ws.Range("C36").Value = "Down arrow key"
$ws.Range("C39").Value = "Ball hits on bar"
$ws.Range("C40").Value = "Ball miss to hit on bar"
$ws.Range("C41").Value = "Game Rematch Screen"
$ws.Range("D41").Value = "game transitions to rematch screen after game loss"
$ws.Range("C42").Value = "Rematch Button Hover"
$ws.Range("C43").Value = "Rematch Button Press"
$ws.Range("D43").Value = "Game restarts"
$ws.Range("B34").Value = "Pong Game"

# Remaining cells of the new block - these reuse strings that already
# exist in the shared string table ("pass", "2", "highlights light grey
# on mouse hover"), so they can be filled in any order.
$ws.Range("A34").Value = 2
$ws.Range("D42").Value = "highlights light grey on mouse hover"
$ws.Range("D44").Value = "highlights light grey on mouse hover"

$ws.Range("E34").Value = "pass"
$ws.Range("E35").Value = "pass"
$ws.Range("E36").Value = "pass"
$ws.Range("E37").Value = "pass"
$ws.Range("E38").Value = "pass"
$ws.Range("E39").Value = "pass"
$ws.Range("E40").Value = "pass"
$ws.Range("E41").Value = "pass"
$ws.Range("E42").Value = "pass"
$ws.Range("E43").Value = "pass"
$ws.Range("E44").Value = "pass"
$ws.Range("E45").Value = "pass"

# ---------------------------------------------------------------------
# 4. The new section header (B34, "Pong Game") uses a smaller Arial font,
#    matching the new font/cellXf added to styles.xml.
# ---------------------------------------------------------------------
$ws.Range("B34").Font.Size = 10
$ws.Range("B34").Font.Name = "Arial"

# ---------------------------------------------------------------------
# 5. Update the view: selected cell moves to C34 and the window scrolls
#    down so row 18 is at the top.
# ---------------------------------------------------------------------
[void]$ws.Range("C34").Select()
$excel.ActiveWindow.ScrollRow = 18
$excel.ActiveWindow.ScrollColumn = 1
